$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column S (19), shifting S:AML right to T:AMM.
# Excel copies the format of the column to the left (R) into the new column,
# matching "Format Same as Left" default insert behavior.
$ws.Columns("S:S").Insert()

# New column header content.
$ws.Range("S1").Value = "Sub brand"

# Selection moves to the newly inserted column's header-adjacent cell.
[void]$ws.Range("S2").Select()
